$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.7310083333333334
$ws.Range("H2").Value = 2.193025
$ws.Range("I2").Value = 0.01673731480740535
$ws.Range("J2").Value = 0.01673731480740535
$ws.Range("M2").Value = 0.13648
$ws.Range("N2").Value = 0.40944
$ws.Range("Q2").Value = 0.09976801733333335
$ws.Range("R2").Value = 0.8979121560000001
$ws.Range("S2").Value = 0.01673731480740535
$ws.Range("T2").Value = 0.01673731480740535

# Row 3
$ws.Range("I3").Value = 0.8536212576586365
$ws.Range("J3").Value = 0.8536212576586365
$ws.Range("M3").Value = 0.13648
$ws.Range("N3").Value = 0.40944
$ws.Range("Q3").Value = 5.0882773856
$ws.Range("R3").Value = 45.7944964704
$ws.Range("S3").Value = 0.8536212576586365
$ws.Range("T3").Value = 0.8536212576586365

# Row 4
$ws.Range("G4").Value = 5.662136666666666
$ws.Range("H4").Value = 16.98641
$ws.Range("I4").Value = 0.129641427533958
$ws.Range("J4").Value = 0.129641427533958
$ws.Range("M4").Value = 0.13648
$ws.Range("N4").Value = 0.40944
$ws.Range("Q4").Value = 0.7727684122666667
$ws.Range("R4").Value = 6.9549157104
$ws.Range("S4").Value = 0.129641427533958
$ws.Range("T4").Value = 0.129641427533958
